$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new test entries (Table_CaptureRatios and Table_UpDownRatios)
# as new rows right after the existing data (rows 123 and 124).
$ws.Range("A123").Value = "Table_CaptureRatios"
$ws.Range("B123").Value = "Test Table_CaptureRatios"
$ws.Range("C123").Value = "Table_CaptureRatios_test"

$ws.Range("A124").Value = "Table_UpDownRatios"
$ws.Range("B124").Value = "Test Table_UpDownRatios"
$ws.Range("C124").Value = "Table_UpDownTatios_test"

# Update the selection to match the post-edit state of the workbook.
$ws.Range("B127").Select()
